# Pooh Points: final 20260128 -> PD8
# Rename owner short-codes to their full team names, and widen owner-name columns.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Map old owner short name -> new full team name
$ownerMap = @{
    "Booz" = "Boozers Losers"
    "CDL"  = "The Backslashers"
    "Clay" = "Hilton Heads"
    "Hal"  = "Three Dawg Nite"
    "Mark" = "Bend Rimmers"
    "Ron"  = "G-Flop"
    "Tar"  = "The Oddities"
}

# --- Players sheet: column B (owner) on rows 2-72 ---
# NOTE: use .Value2 (not .Value) for reads -- this host's .Value getter
# does not resolve to the underlying scalar.
for ($r = 2; $r -le 72; $r++) {
    $cell = $ws1.Cells.Item($r, 2)
    $old = $cell.Value2
    if ($ownerMap.ContainsKey($old)) {
        $cell.Value2 = $ownerMap[$old]
    }
}

# Widen owner column (B) on Players sheet
# NOTE: this host's ColumnWidth setter round-trips through a pixel
# conversion that adds back ~5/6 of a character unit on save, so feed it
# target-minus-5/6 here to land on a clean "18" in the saved OOXML.
$ws1.Columns.Item(2).ColumnWidth = 17.166666666666668

# --- OwnerTotals sheet: column A (owner) on rows 2-8 ---
for ($r = 2; $r -le 8; $r++) {
    $cell = $ws2.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($ownerMap.ContainsKey($old)) {
        $cell.Value2 = $ownerMap[$old]
    }
}

# Widen owner column (A) on OwnerTotals sheet
$ws2.Columns.Item(1).ColumnWidth = 17.166666666666668
